$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.914.94'
$ws.Range("E2").Value = '  +0.78%  '

$ws.Range("D3").Value = '3.392.79'
$ws.Range("E3").Value = '  +0.36%  '

$ws.Range("E4").Value = '  +0.13%  '

$ws.Range("D5").Value = '''579.94'
$ws.Range("E5").Value = '  +0.54%  '

$ws.Range("D6").Value = '''138.16'
$ws.Range("E6").Value = '  +2.28%  '

$ws.Range("E7").Value = '  -0.01%  '

$ws.Range("D8").Value = '3.395.02'
$ws.Range("E8").Value = '  +0.44%  '

$ws.Range("E9").Value = '  -0.76%  '

$ws.Range("D10").Value = '''7.51'
$ws.Range("E10").Value = '  -1.18%  '

$ws.Range("E11").Value = '  +3.04%  '

$ws.Range("E12").Value = '  +0.87%  '

$ws.Range("D13").Value = '3.975.48'
$ws.Range("E13").Value = '  +0.63%  '

$ws.Range("E14").Value = '  +1.39%  '

$ws.Range("D15").Value = '''0.0000178'
$ws.Range("E15").Value = '  +1.46%  '

$ws.Range("D16").Value = '3.396.73'
$ws.Range("E16").Value = '  +0.60%  '

$ws.Range("D17").Value = '''25.39'
$ws.Range("E17").Value = '  +0.58%  '

$ws.Range("D18").Value = '62.067.57'
$ws.Range("E18").Value = '  +1.13%  '

$ws.Range("D19").Value = '''14.15'
$ws.Range("E19").Value = '  +0.74%  '

$ws.Range("B20").Value = 'Polkadot'
$ws.Range("C20").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D20").Value = '''5.81'
$ws.Range("E20").Value = '  -0.20%  '

$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").Value = '''9.44'
$ws.Range("E21").Value = '  +0.88%  '

$ws.Range("D22").Value = '''395.00'
$ws.Range("E22").Value = '  +3.89%  '

$ws.Range("D23").Value = '''0.565'
$ws.Range("E23").Value = '  -0.60%  '

$ws.Range("B24").Value = 'WrappedeETH'
$ws.Range("C24").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D24").Value = '3.544.77'
$ws.Range("E24").Value = '  +0.99%  '

$ws.Range("B25").Value = 'PEPE'
$ws.Range("C25").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D25").Value = '''0.0000129'
$ws.Range("E25").Value = '  +8.19%  '

$ws.Range("D26").Value = '''0.999'
$ws.Range("E26").Value = '  -0.12%  '

$ws.Range("D27").Value = '''71.59'
$ws.Range("E27").Value = '  +1.18%  '

$ws.Range("D28").Value = '''1.65'
$ws.Range("E28").Value = '  -2.57%  '

$ws.Range("D29").Value = '''7.64'
$ws.Range("E29").Value = '  -2.12%  '

$ws.Range("D30").Value = '''1.00'
$ws.Range("E30").Value = '  +0.19%  '

$ws.Range("E31").Value = '  +3.41%  '

$ws.Range("D32").Value = '''8.22'
$ws.Range("E32").Value = '  +0.63%  '

$ws.Range("D33").Value = '''2.18'
$ws.Range("E33").Value = '  +0.23%  '

$ws.Range("E34").Value = '  -0.01%  '

$ws.Range("D35").Value = '''23.48'
$ws.Range("E35").Value = '  +0.20%  '

$ws.Range("D36").Value = '3.427.54'
$ws.Range("E36").Value = '  +0.53%  '

$ws.Range("D37").Value = '''5.39'
$ws.Range("E37").Value = '  -3.69%  '

$ws.Range("D38").Value = '''1.59'
$ws.Range("E38").Value = '  +2.19%  '

$ws.Range("D39").Value = '''6.90'
$ws.Range("E39").Value = '  -1.42%  '

$ws.Range("D40").Value = '''164.90'
$ws.Range("E40").Value = '  +1.22%  '

$ws.Range("D41").Value = '''0.0788'
$ws.Range("E41").Value = '  +0.24%  '

$ws.Range("D42").Value = '''1.77'
$ws.Range("E42").Value = '  +8.77%  '

$ws.Range("D43").Value = '''1.25'
$ws.Range("E43").Value = '  +1.43%  '

$ws.Range("D44").Value = '''0.787'
$ws.Range("E44").Value = '  +3.70%  '

$ws.Range("E45").Value = '  +0.23%  '

$ws.Range("D46").Value = '''4.43'
$ws.Range("E46").Value = '  -0.35%  '

$ws.Range("D47").Value = '''24.97'
$ws.Range("E47").Value = '  +6.62%  '

$ws.Range("D48").Value = '''41.30'
$ws.Range("E48").Value = '  -0.97%  '

$ws.Range("D49").Value = '''6.89'
$ws.Range("E49").Value = '  -0.87%  '

$ws.Range("D50").Value = '''23.04'
$ws.Range("E50").Value = '  -1.51%  '

$ws.Range("D51").Value = '2.337.20'
$ws.Range("E51").Value = '  +6.34%  '
